$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (row order + updated total_registros values) for A2:B15
$data = @(
    @("ALBIRENA GARCIA ANGEELO ALONSO", 146),
    @("MANUEL LEUNARDO PRADO BAILON", 133),
    @("ALAMA NIMA CLARITZA MABEL", 120),
    @("URRIOLA ARISMENDIZ INGRID MARYURI", 110),
    @("CORDOVA CARMEN ANGIE NATALLY", 108),
    @("AGURTO ORDINOLA LISBET JAQUELIN", 104),
    @("CARREÑO PALACIOS KATHERINE DE LOS MILAGROS", 99),
    @("MARYURI OJEDA VALLE", 96),
    @("VEGA ROBLEDO FERNANDO ERNESTO", 96),
    @("JUAREZ CARMEN PIERRE ALEXANDER", 93),
    @("ATOCHE PALACIOS LUIS ANGEL", 91),
    @("ROMAN GALECIO MARITZA DEL PILAR", 85),
    @("BERNAOLA CARMEN ZUMIKO YASHURY", 76),
    @("RUIDIAS FRIAS MELISSA VICTORIA", 73)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# Remove the old trailing row 16 (61097774 / 1) entirely, shrinking the used range to A1:B15
$ws.Rows.Item(16).Delete()
